$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells from camelCase to snake_case. The column order
# (A=sepalLength/B=sepalWidth/C=petalLength/D=petalWidth/E=irisClass) stays
# the same -- only the header text changes. The trailing portion of each
# new header (starting at the underscore) is re-bolded, matching the rich
# text runs the author applied in Excel.

$ws.Range("A1").Value = "sepal_length"
$ws.Range("A1").Characters(6, 7).Font.Bold = $true

$ws.Range("B1").Value = "sepal_width"
$ws.Range("B1").Characters(6, 6).Font.Bold = $true

$ws.Range("C1").Value = "petal_length"
$ws.Range("C1").Characters(6, 7).Font.Bold = $true

$ws.Range("D1").Value = "petal_width"
$ws.Range("D1").Characters(6, 6).Font.Bold = $true

$ws.Range("E1").Value = "iris_class"
$ws.Range("E1").Characters(5, 6).Font.Bold = $true

# Resize columns to fit the new header text, as seen in the saved workbook.
$ws.Columns.AutoFit() | Out-Null
